$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B (Black carbon aerosol)
$ws.Range("B2").Value  = "+7.76 ± 0.00"
$ws.Range("B3").Value  = "+7.76 ± 0.00"
$ws.Range("B5").Value  = "-7.77 ± 0.00"
$ws.Range("B7").Value  = "-1.02 ± 0.00"
$ws.Range("B9").Value  = "-6.31 ± 0.00"
$ws.Range("B10").Value = "-0.42 ± 0.00"
$ws.Range("B11").Value = "+0.11 ± 0.00"
$ws.Range("B12").Value = "+5.19 ± 0.01"

# Column C (in pure BC)
$ws.Range("C2").Value  = "+7.76 ± 0.00"
$ws.Range("C3").Value  = "+7.76 ± 0.00"
$ws.Range("C5").Value  = "-7.76 ± 0.00"
$ws.Range("C6").Value  = "-1.83 ± 0.00"
$ws.Range("C9").Value  = "-5.56 ± 0.00"
$ws.Range("C10").Value = "-0.36 ± 0.00"
$ws.Range("C11").Value = "+0.09 ± 0.00"
$ws.Range("C12").Value = "+4.32 ± 0.01"

# Column D (in MBS)
$ws.Range("D2").Value  = "+1.83 ± 0.00"
$ws.Range("D4").Value  = "+1.83 ± 0.00"
$ws.Range("D5").Value  = "-1.84 ± 0.00"
$ws.Range("D7").Value  = "-1.02 ± 0.00"
$ws.Range("D9").Value  = "-0.74 ± 0.00"
$ws.Range("D12").Value = "+3.69 ± 0.01"
